# Remove the first data row (chr6 / 26555436 / 26555459 / ...).
# All subsequent rows shift up by one, so row 2 becomes the new row 1,
# row 9 becomes the new row 8, and the sheet's used range shrinks from
# A1:S9 to A1:S8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Delete()
